$wb = $excel.ActiveWorkbook
$source = $wb.Worksheets.Item("strategy_id-5009")

$sheetNames = @("strategy_id-5011", "strategy_id-5012", "strategy_id-5013", "strategy_id-5014", "strategy_id-5015")

foreach ($name in $sheetNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $source.Copy([System.Reflection.Missing]::Value, $lastSheet)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $name
}

$ws = $wb.Worksheets.Item("strategy_id-5011")
$ws.Cells.Item(2, 21).Value = 0.9829063213546697
$ws.Cells.Item(2, 22).Value = 0.9658126427093395
$ws.Cells.Item(2, 23).Value = 0.9487189640640092
$ws.Cells.Item(2, 24).Value = 0.931625285418679
$ws.Cells.Item(2, 25).Value = 0.9145316067733488
$ws.Cells.Item(2, 26).Value = 0.8974379281280185
$ws.Cells.Item(2, 27).Value = 0.8803442494826883
$ws.Cells.Item(2, 28).Value = 0.863250570837358
$ws.Cells.Item(2, 29).Value = 0.8461568921920278
$ws.Cells.Item(2, 30).Value = 0.8290632135466975
$ws.Cells.Item(2, 31).Value = 0.8119695349013674
$ws.Cells.Item(2, 32).Value = 0.7948758562560371
$ws.Cells.Item(2, 33).Value = 0.7777821776107068
$ws.Cells.Item(2, 34).Value = 0.7606884989653766
$ws.Cells.Item(2, 35).Value = 0.7435948203200462
$ws.Cells.Item(2, 36).Value = 0.726501141674716
$ws.Cells.Item(2, 37).Value = 0.7094074630293857
$ws.Cells.Item(2, 38).Value = 0.6923137843840557
$ws.Cells.Item(2, 39).Value = 0.6752201057387253
$ws.Cells.Item(2, 40).Value = 0.6581264270933951
$ws.Cells.Item(2, 41).Value = 0.6410327484480649
$ws.Cells.Item(2, 42).Value = 0.6239390698027346
$ws.Cells.Item(2, 43).Value = 0.6068453911574043
$ws.Cells.Item(2, 44).Value = 0.5897517125120741
$ws.Cells.Item(2, 45).Value = 0.5726580338667439

$ws = $wb.Worksheets.Item("strategy_id-5012")
$ws.Cells.Item(2, 21).Value = 0.9862304763364699
$ws.Cells.Item(2, 22).Value = 0.9724609526729397
$ws.Cells.Item(2, 23).Value = 0.9586914290094095
$ws.Cells.Item(2, 24).Value = 0.9449219053458794
$ws.Cells.Item(2, 25).Value = 0.9311523816823493
$ws.Cells.Item(2, 26).Value = 0.9173828580188191
$ws.Cells.Item(2, 27).Value = 0.903613334355289
$ws.Cells.Item(2, 28).Value = 0.8898438106917588
$ws.Cells.Item(2, 29).Value = 0.8760742870282288
$ws.Cells.Item(2, 30).Value = 0.8623047633646985
$ws.Cells.Item(2, 31).Value = 0.8485352397011684
$ws.Cells.Item(2, 32).Value = 0.8347657160376383
$ws.Cells.Item(2, 33).Value = 0.820996192374108
$ws.Cells.Item(2, 34).Value = 0.8072266687105779
$ws.Cells.Item(2, 35).Value = 0.7934571450470478
$ws.Cells.Item(2, 36).Value = 0.7796876213835175
$ws.Cells.Item(2, 37).Value = 0.7659180977199874
$ws.Cells.Item(2, 38).Value = 0.7521485740564573
$ws.Cells.Item(2, 39).Value = 0.7383790503929272
$ws.Cells.Item(2, 40).Value = 0.724609526729397
$ws.Cells.Item(2, 41).Value = 0.7108400030658669
$ws.Cells.Item(2, 42).Value = 0.6970704794023368
$ws.Cells.Item(2, 43).Value = 0.6833009557388066
$ws.Cells.Item(2, 44).Value = 0.6695314320752765
$ws.Cells.Item(2, 45).Value = 0.6557619084117463

$ws = $wb.Worksheets.Item("strategy_id-5013")
$ws.Cells.Item(2, 21).Value = 0.9831702487644064
$ws.Cells.Item(2, 22).Value = 0.9663404975288129
$ws.Cells.Item(2, 23).Value = 0.9495107462932192
$ws.Cells.Item(2, 24).Value = 0.9326809950576255
$ws.Cells.Item(2, 25).Value = 0.915851243822032
$ws.Cells.Item(2, 26).Value = 0.8990214925864384
$ws.Cells.Item(2, 27).Value = 0.8821917413508448
$ws.Cells.Item(2, 28).Value = 0.8653619901152512
$ws.Cells.Item(2, 29).Value = 0.8485322388796576
$ws.Cells.Item(2, 30).Value = 0.831702487644064
$ws.Cells.Item(2, 31).Value = 0.8148727364084704
$ws.Cells.Item(2, 32).Value = 0.7980429851728768
$ws.Cells.Item(2, 33).Value = 0.7812132339372831
$ws.Cells.Item(2, 34).Value = 0.7643834827016895
$ws.Cells.Item(2, 35).Value = 0.7475537314660958
$ws.Cells.Item(2, 36).Value = 0.7307239802305022
$ws.Cells.Item(2, 37).Value = 0.7138942289949086
$ws.Cells.Item(2, 38).Value = 0.6970644777593151
$ws.Cells.Item(2, 39).Value = 0.6802347265237214
$ws.Cells.Item(2, 40).Value = 0.6634049752881279
$ws.Cells.Item(2, 41).Value = 0.6465752240525343
$ws.Cells.Item(2, 42).Value = 0.6297454728169407
$ws.Cells.Item(2, 43).Value = 0.6129157215813471
$ws.Cells.Item(2, 44).Value = 0.5960859703457535
$ws.Cells.Item(2, 45).Value = 0.5792562191101598

$ws = $wb.Worksheets.Item("strategy_id-5014")
$ws.Cells.Item(2, 21).Value = 0.9888585559026216
$ws.Cells.Item(2, 22).Value = 0.9777171118052433
$ws.Cells.Item(2, 23).Value = 0.966575667707865
$ws.Cells.Item(2, 24).Value = 0.9554342236104866
$ws.Cells.Item(2, 25).Value = 0.9442927795131083
$ws.Cells.Item(2, 26).Value = 0.9331513354157299
$ws.Cells.Item(2, 27).Value = 0.9220098913183516
$ws.Cells.Item(2, 28).Value = 0.9108684472209733
$ws.Cells.Item(2, 29).Value = 0.899727003123595
$ws.Cells.Item(2, 30).Value = 0.8885855590262166
$ws.Cells.Item(2, 31).Value = 0.8774441149288383
$ws.Cells.Item(2, 32).Value = 0.86630267083146
$ws.Cells.Item(2, 33).Value = 0.8551612267340816
$ws.Cells.Item(2, 34).Value = 0.8440197826367033
$ws.Cells.Item(2, 35).Value = 0.8328783385393249
$ws.Cells.Item(2, 36).Value = 0.8217368944419465
$ws.Cells.Item(2, 37).Value = 0.8105954503445683
$ws.Cells.Item(2, 38).Value = 0.79945400624719
$ws.Cells.Item(2, 39).Value = 0.7883125621498116
$ws.Cells.Item(2, 40).Value = 0.7771711180524332
$ws.Cells.Item(2, 41).Value = 0.766029673955055
$ws.Cells.Item(2, 42).Value = 0.7548882298576767
$ws.Cells.Item(2, 43).Value = 0.7437467857602983
$ws.Cells.Item(2, 44).Value = 0.7326053416629199
$ws.Cells.Item(2, 45).Value = 0.7214638975655416

$ws = $wb.Worksheets.Item("strategy_id-5015")
$ws.Cells.Item(2, 21).Value = 0.9813124864730004
$ws.Cells.Item(2, 22).Value = 0.9626249729460008
$ws.Cells.Item(2, 23).Value = 0.9439374594190011
$ws.Cells.Item(2, 24).Value = 0.9252499458920015
$ws.Cells.Item(2, 25).Value = 0.906562432365002
$ws.Cells.Item(2, 26).Value = 0.8878749188380024
$ws.Cells.Item(2, 27).Value = 0.8691874053110027
$ws.Cells.Item(2, 28).Value = 0.8504998917840031
$ws.Cells.Item(2, 29).Value = 0.8318123782570036
$ws.Cells.Item(2, 30).Value = 0.813124864730004
$ws.Cells.Item(2, 31).Value = 0.7944373512030043
$ws.Cells.Item(2, 32).Value = 0.7757498376760047
$ws.Cells.Item(2, 33).Value = 0.7570623241490051
$ws.Cells.Item(2, 34).Value = 0.7383748106220054
$ws.Cells.Item(2, 35).Value = 0.7196872970950059
$ws.Cells.Item(2, 36).Value = 0.7009997835680062
$ws.Cells.Item(2, 37).Value = 0.6823122700410066
$ws.Cells.Item(2, 38).Value = 0.663624756514007
$ws.Cells.Item(2, 39).Value = 0.6449372429870075
$ws.Cells.Item(2, 40).Value = 0.6262497294600078
$ws.Cells.Item(2, 41).Value = 0.6075622159330083
$ws.Cells.Item(2, 42).Value = 0.5888747024060086
$ws.Cells.Item(2, 43).Value = 0.5701871888790091
$ws.Cells.Item(2, 44).Value = 0.5514996753520094
$ws.Cells.Item(2, 45).Value = 0.5328121618250098
